# Apply MassWateR site-build edit: add two red "updated" notes to the
# Instructions sheet (column C, rows 1-2) documenting when the Meta
# template / sample data were last refreshed.

$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Meta")
$wsInstructions = $wb.Worksheets.Item("Instructions")

$wsInstructions.Range("C1").Value = "Template updated 1/17/23"
$wsInstructions.Range("C2").Value = "Samples updated 1/8/23"

$noteRange = $wsInstructions.Range("C1:C2")
$noteRange.Font.Color = 255

$wsMeta.Range("B12").Select()
$wsInstructions.Range("C3").Select()
$wsMeta.Activate()
